$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the row for "nafo informado" (row 320). Deleting it shifts all
# subsequent rows up by one, matching the new dimension A1:C562.
$ws.Rows.Item(320).Delete()
